$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lpl"
$ws.Range("C2").Value = "Vldlr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 46.48074800000001
$ws.Range("H2").Value = 139.442244
$ws.Range("I2").Value = 0.1473944418036112
$ws.Range("J2").Value = 0.1473944418036112
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1112246666666667
$ws.Range("N2").Value = 0.333674
$ws.Range("O2").Value = 0.008007360146089436
$ws.Range("P2").Value = 0.008007360146089434
$ws.Range("Q2").Value = 5.169805702717334
$ws.Range("R2").Value = 46.52825132445601
$ws.Range("S2").Value = 0.001180240379053335
$ws.Range("T2").Value = 0.001180240379053334

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lpl"
$ws.Range("C3").Value = "Vldlr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 46.48074800000001
$ws.Range("H3").Value = 139.442244
$ws.Range("I3").Value = 0.1473944418036112
$ws.Range("J3").Value = 0.1473944418036112
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.825228
$ws.Range("N3").Value = 35.475684
$ws.Range("O3").Value = 0.851329675722
$ws.Range("P3").Value = 0.8513296757219999
$ws.Range("Q3").Value = 549.6454427105441
$ws.Range("R3").Value = 4946.808984394896
$ws.Range("S3").Value = 0.1254812623438935
$ws.Range("T3").Value = 0.1254812623438935

# Row 4: ECs -> M2
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lpl"
$ws.Range("C4").Value = "Vldlr"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 46.48074800000001
$ws.Range("H4").Value = 139.442244
$ws.Range("I4").Value = 0.1473944418036112
$ws.Range("J4").Value = 0.1473944418036112
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03793733333333333
$ws.Range("N4").Value = 0.113812
$ws.Range("O4").Value = 0.002731209722503793
$ws.Range("P4").Value = 0.002731209722503793
$ws.Range("Q4").Value = 1.763355630458667
$ws.Range("R4").Value = 15.870200674128
$ws.Range("S4").Value = 0.0004025651324970423
$ws.Range("T4").Value = 0.0004025651324970423

# Row 5: ECs -> sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Lpl"
$ws.Range("C5").Value = "Vldlr"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 46.48074800000001
$ws.Range("H5").Value = 139.442244
$ws.Range("I5").Value = 0.1473944418036112
$ws.Range("J5").Value = 0.1473944418036112
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.915914
$ws.Range("N5").Value = 5.747742000000001
$ws.Range("O5").Value = 0.1379317544094067
$ws.Range("P5").Value = 0.1379317544094067
$ws.Range("Q5").Value = 89.05311582367202
$ws.Range("R5").Value = 801.4780424130482
$ws.Range("S5").Value = 0.02033037394816729
$ws.Range("T5").Value = 0.02033037394816729

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lpl"
$ws.Range("C6").Value = "Vldlr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 192.830597
$ws.Range("H6").Value = 578.4917909999999
$ws.Range("I6").Value = 0.6114823756165045
$ws.Range("J6").Value = 0.6114823756165044
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1112246666666667
$ws.Range("N6").Value = 0.333674
$ws.Range("O6").Value = 0.008007360146089436
$ws.Range("P6").Value = 0.008007360146089434
$ws.Range("Q6").Value = 21.44751887445933
$ws.Range("R6").Value = 193.027669870134
$ws.Range("S6").Value = 0.004896359604547689
$ws.Range("T6").Value = 0.004896359604547687

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lpl"
$ws.Range("C7").Value = "Vldlr"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 192.830597
$ws.Range("H7").Value = 578.4917909999999
$ws.Range("I7").Value = 0.6114823756165045
$ws.Range("J7").Value = 0.6114823756165044
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 11.825228
$ws.Range("N7").Value = 35.475684
$ws.Range("O7").Value = 0.851329675722
$ws.Range("P7").Value = 0.8513296757219999
$ws.Range("Q7").Value = 2280.265774901116
$ws.Range("R7").Value = 20522.39197411004
$ws.Range("S7").Value = 0.5205730925433171
$ws.Range("T7").Value = 0.5205730925433169

# Row 8: FAPs -> M2
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Lpl"
$ws.Range("C8").Value = "Vldlr"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 192.830597
$ws.Range("H8").Value = 578.4917909999999
$ws.Range("I8").Value = 0.6114823756165045
$ws.Range("J8").Value = 0.6114823756165044
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.03793733333333333
$ws.Range("N8").Value = 0.113812
$ws.Range("O8").Value = 0.002731209722503793
$ws.Range("P8").Value = 0.002731209722503793
$ws.Range("Q8").Value = 7.315478635254665
$ws.Range("R8").Value = 65.83930771729199
$ws.Range("S8").Value = 0.001670086609423513
$ws.Range("T8").Value = 0.001670086609423513

# Row 9: FAPs -> sCs
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Lpl"
$ws.Range("C9").Value = "Vldlr"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 192.830597
$ws.Range("H9").Value = 578.4917909999999
$ws.Range("I9").Value = 0.6114823756165045
$ws.Range("J9").Value = 0.6114823756165044
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.915914
$ws.Range("N9").Value = 5.747742000000001
$ws.Range("O9").Value = 0.1379317544094067
$ws.Range("P9").Value = 0.1379317544094067
$ws.Range("Q9").Value = 369.446840420658
$ws.Range("R9").Value = 3325.021563785922
$ws.Range("S9").Value = 0.0843428368592163
$ws.Range("T9").Value = 0.08434283685921629

# Row 10: M2 -> ECs
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Lpl"
$ws.Range("C10").Value = "Vldlr"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 60.030993
$ws.Range("H10").Value = 180.092979
$ws.Range("I10").Value = 0.1903634318482028
$ws.Range("J10").Value = 0.1903634318482028
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1112246666666667
$ws.Range("N10").Value = 0.333674
$ws.Range("O10").Value = 0.008007360146089436
$ws.Range("P10").Value = 0.008007360146089434
$ws.Range("Q10").Value = 6.676927186094001
$ws.Range("R10").Value = 60.09234467484601
$ws.Range("S10").Value = 0.001524308557454112
$ws.Range("T10").Value = 0.001524308557454111

# Row 11: M2 -> FAPs
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Lpl"
$ws.Range("C11").Value = "Vldlr"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 60.030993
$ws.Range("H11").Value = 180.092979
$ws.Range("I11").Value = 0.1903634318482028
$ws.Range("J11").Value = 0.1903634318482028
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 11.825228
$ws.Range("N11").Value = 35.475684
$ws.Range("O11").Value = 0.851329675722
$ws.Range("P11").Value = 0.8513296757219999
$ws.Range("Q11").Value = 709.880179291404
$ws.Range("R11").Value = 6388.921613622637
$ws.Range("S11").Value = 0.1620620387046575
$ws.Range("T11").Value = 0.1620620387046575

# Row 12: M2 -> M2
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Lpl"
$ws.Range("C12").Value = "Vldlr"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 60.030993
$ws.Range("H12").Value = 180.092979
$ws.Range("I12").Value = 0.1903634318482028
$ws.Range("J12").Value = 0.1903634318482028
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.03793733333333333
$ws.Range("N12").Value = 0.113812
$ws.Range("O12").Value = 0.002731209722503793
$ws.Range("P12").Value = 0.002731209722503793
$ws.Range("Q12").Value = 2.277415791772
$ws.Range("R12").Value = 20.496742125948
$ws.Range("S12").Value = 0.0005199224558729998
$ws.Range("T12").Value = 0.0005199224558729996

# Row 13: M2 -> sCs
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Lpl"
$ws.Range("C13").Value = "Vldlr"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 60.030993
$ws.Range("H13").Value = 180.092979
$ws.Range("I13").Value = 0.1903634318482028
$ws.Range("J13").Value = 0.1903634318482028
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.915914
$ws.Range("N13").Value = 5.747742000000001
$ws.Range("O13").Value = 0.1379317544094067
$ws.Range("P13").Value = 0.1379317544094067
$ws.Range("Q13").Value = 115.014219922602
$ws.Range("R13").Value = 1035.127979303418
$ws.Range("S13").Value = 0.02625716213021814
$ws.Range("T13").Value = 0.02625716213021814

# Row 14: sCs -> ECs
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Lpl"
$ws.Range("C14").Value = "Vldlr"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 16.00705666666667
$ws.Range("H14").Value = 48.02117
$ws.Range("I14").Value = 0.05075975073168155
$ws.Range("J14").Value = 0.05075975073168155
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.1112246666666667
$ws.Range("N14").Value = 0.333674
$ws.Range("O14").Value = 0.008007360146089436
$ws.Range("P14").Value = 0.008007360146089434
$ws.Range("Q14").Value = 1.780379542064445
$ws.Range("R14").Value = 16.02341587858
$ws.Range("S14").Value = 0.0004064516050343009
$ws.Range("T14").Value = 0.0004064516050343008

# Row 15: sCs -> FAPs
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Lpl"
$ws.Range("C15").Value = "Vldlr"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 16.00705666666667
$ws.Range("H15").Value = 48.02117
$ws.Range("I15").Value = 0.05075975073168155
$ws.Range("J15").Value = 0.05075975073168155
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 11.825228
$ws.Range("N15").Value = 35.475684
$ws.Range("O15").Value = 0.851329675722
$ws.Range("P15").Value = 0.8513296757219999
$ws.Range("Q15").Value = 189.2870946922534
$ws.Range("R15").Value = 1703.58385223028
$ws.Range("S15").Value = 0.04321328213013201
$ws.Range("T15").Value = 0.043213282130132

# Row 16: sCs -> M2
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Lpl"
$ws.Range("C16").Value = "Vldlr"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 16.00705666666667
$ws.Range("H16").Value = 48.02117
$ws.Range("I16").Value = 0.05075975073168155
$ws.Range("J16").Value = 0.05075975073168155
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.03793733333333333
$ws.Range("N16").Value = 0.113812
$ws.Range("O16").Value = 0.002731209722503793
$ws.Range("P16").Value = 0.002731209722503793
$ws.Range("Q16").Value = 0.6072650444488888
$ws.Range("R16").Value = 5.46538540004
$ws.Range("S16").Value = 0.0001386355247102377
$ws.Range("T16").Value = 0.0001386355247102377

# Row 17: sCs -> sCs
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Lpl"
$ws.Range("C17").Value = "Vldlr"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 16.00705666666667
$ws.Range("H17").Value = 48.02117
$ws.Range("I17").Value = 0.05075975073168155
$ws.Range("J17").Value = 0.05075975073168155
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.915914
$ws.Range("N17").Value = 5.747742000000001
$ws.Range("O17").Value = 0.1379317544094067
$ws.Range("P17").Value = 0.1379317544094067
$ws.Range("Q17").Value = 30.66814396646
$ws.Range("R17").Value = 276.01329569814
$ws.Range("S17").Value = 0.007001381471805003
$ws.Range("T17").Value = 0.007001381471805002
